# Update cryptocurrency price/volume data per the latest GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.035.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.489.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.05%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.98%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -2.92%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.488.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.05%  "

$ws.Range("E10").Value = "  -2.28%  "

$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("E12").Value = "  -4.56%  "

$ws.Range("E13").Value = "  -2.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.939.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.859.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.490.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.69%  "

$ws.Range("E22").Value = "  -3.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.58%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "68.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.613.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0958"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "520.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.70%  "

$ws.Range("E34").Value = "  -4.36%  "

$ws.Range("E35").Value = "  -3.41%  "

$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("E38").Value = "  -3.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "

$ws.Range("E41").Value = "  -3.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.46%  "

$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "145.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.35%  "

$ws.Range("E47").Value = "  -4.79%  "

$ws.Range("E48").Value = "  -3.48%  "

$ws.Range("E49").Value = "  -9.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0748"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.68%  "
